# Daily attendance processing - 2025-10-16 20:31:29
#
# The "Recorded By" column (G) lists the users who touched each attendance
# session, as a comma-separated string. The attendance sync re-orders this
# list so that the first contributor is rotated to the end (first-in,
# first-out) for every multi-author row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$col = 7  # column G = "Recorded By"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $text = $cell.Text

    if ($text -and $text.Contains(",")) {
        $parts = $text -split ",\s*"
        $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
        $cell.Value = $rotated
    }
}
